# dbow frequent word sampling & slightly better results
#
# - Change E12/E13 from numeric 10 to the shared "-" string (frequent word
#   sampling wasn't actually used for those dbow runs).
# - Append 4 new dbow rows (14-17) covering the w2v frequent-word-sampling
#   thresholds 1e-2 / 1e-3 / 1e-4 / 1e-5, each at window size "-" (no
#   window size column data for dbow).
# - Move the active-cell selection to K15 to match the new data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doc vecs")

# --- fix up E12 / E13: "window size" becomes "-" (was 10) -----------------
$ws.Range("E12").Value = "-"
$ws.Range("E13").Value = "-"

# --- new row 14: dbow, w2v 1e-2 -------------------------------------------
$ws.Range("A14").Value = "dbow"
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = "w2v 1e-2"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = 0.001
$ws.Range("G14").Value = 0.001
$ws.Range("H14").Value = 20
$ws.Range("I14").Value = 10
$ws.Range("I14").NumberFormat = "0"
$ws.Range("J14").Value = 0.024305555555555556
$ws.Range("J14").NumberFormat = "h:mm"
$ws.Range("K14").Value = 0.0055555555555555558
$ws.Range("K14").NumberFormat = "h:mm"
$ws.Range("L14").Value = 0.11935999999999999

# --- new row 15: dbow, w2v 1e-3 -------------------------------------------
$ws.Range("A15").Style = "Normal"
$ws.Range("A15").Value = "dbow"
$ws.Range("B15").Value = 100
$ws.Range("C15").Value = "w2v 1e-3"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = 0.001
$ws.Range("G15").Value = 0.001
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 10
$ws.Range("I15").NumberFormat = "0"
$ws.Range("J15").Value = 0.019444444444444445
$ws.Range("J15").NumberFormat = "h:mm"
$ws.Range("K15").Value = 0.0041666666666666666
$ws.Range("K15").NumberFormat = "h:mm"
$ws.Range("L15").Value = 0.11744

# --- new row 16: dbow, w2v 1e-4 -------------------------------------------
$ws.Range("A16").Value = "dbow"
$ws.Range("B16").Value = 100
$ws.Range("C16").Value = "w2v 1e-4"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = 0.001
$ws.Range("G16").Value = 0.001
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 10
$ws.Range("I16").NumberFormat = "0"
$ws.Range("J16").Value = 0.014583333333333332
$ws.Range("J16").NumberFormat = "h:mm"
$ws.Range("K16").Value = 0.0027777777777777779
$ws.Range("K16").NumberFormat = "h:mm"
$ws.Range("L16").Value = 0.12923999999999999

# --- new row 17: dbow, w2v 1e-5 -------------------------------------------
$ws.Range("A17").Value = "dbow"
$ws.Range("B17").Value = 100
$ws.Range("C17").Value = "w2v 1e-5"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "-"
$ws.Range("F17").Value = 0.001
$ws.Range("G17").Value = 0.001
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 10
$ws.Range("I17").NumberFormat = "0"
$ws.Range("J17").Value = 0.0083333333333333332
$ws.Range("J17").NumberFormat = "h:mm"
$ws.Range("K17").Value = 0.0020833333333333333
$ws.Range("K17").NumberFormat = "h:mm"
$ws.Range("L17").Value = 0.23119999999999999

# --- selection follows the newly entered data ------------------------------
$ws.Range("K15").Select()
